$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Registrar elogio aos servidores"
$ws.Range("B3").Value = "https://www.ms.gov.br/comunicacao-e-transparencia/teste-4189"
$ws.Range("C3").Value = "Seção 'Outras Informações' não encontrada"
